$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '26.643.49'
$ws.Range("E2").Value = '  +1.53%  '

# Row 3
$ws.Range("D3").Value = '1.628.70'
$ws.Range("E3").Value = '  +1.72%  '

# Row 4
$ws.Range("E4").Value = '  -0.06%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '213.29'
$ws.Range("E5").Value = '  +0.46%  '

# Row 6
$ws.Range("E6").Value = '  -0.04%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.487'
$ws.Range("E7").Value = '  +0.81%  '

# Row 8
$ws.Range("E8").Value = '  +1.07%  '

# Row 9
$ws.Range("E9").Value = '  +1.02%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.01'
$ws.Range("E10").Value = '  +5.00%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0831'
$ws.Range("E11").Value = '  +2.39%  '

# Row 13
$ws.Range("D13").Value = '1.602.66'
$ws.Range("E13").Value = '  +0.06%  '

# Row 14
$ws.Range("E14").Value = '  +0.55%  '

# Row 15
$ws.Range("E15").Value = '  +2.22%  '

# Row 16
$ws.Range("D16").Value = '26.619.57'
$ws.Range("E16").Value = '  +1.52%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '63.06'
$ws.Range("E17").Value = '  +2.97%  '

# Row 18
$ws.Range("D18").Value = '0.0₃0732'
$ws.Range("E18").Value = '  +0.39%  '

# Row 19
$ws.Range("B19").Value = 'BitcoinCash'
$ws.Range("C19").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '208.49'
$ws.Range("E19").Value = '  +3.13%  '

# Row 20
$ws.Range("B20").Value = 'Dai'
$ws.Range("C20").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.00'
$ws.Range("E20").Value = '  -0.20%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.31'
$ws.Range("E21").Value = '  +0.90%  '

# Row 22
$ws.Range("E22").Value = '  +1.70%  '

# Row 23
$ws.Range("E23").Value = '  +1.65%  '

# Row 24
$ws.Range("E24").Value = '  -1.68%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '145.28'
$ws.Range("E25").Value = '  +0.72%  '

# Row 26
$ws.Range("E26").Value = '  -0.13%  '

# Row 27
$ws.Range("E27").Value = '  -0.94%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '15.42'
$ws.Range("E28").Value = '  +1.69%  '

# Row 29
$ws.Range("E29").Value = '  +1.54%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0521'
$ws.Range("E30").Value = '  +6.87%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.17'
$ws.Range("E31").Value = '  +0.62%  '

# Row 32
$ws.Range("E32").Value = '  +1.89%  '

# Row 33
$ws.Range("E33").Value = '  +0.81%  '

# Row 34
$ws.Range("E34").Value = '  +1.66%  '

# Row 35
$ws.Range("E35").Value = '  -0.27%  '

# Row 36
$ws.Range("D36").Value = '1.162.47'
$ws.Range("E36").Value = '  +0.49%  '

# Row 37
$ws.Range("E37").Value = '  -0.35%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.806'
$ws.Range("E38").Value = '  +1.95%  '

# Row 39
$ws.Range("E39").Value = '  -0.12%  '

# Row 40
$ws.Range("E40").Value = '  -0.07%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.501'
$ws.Range("E41").Value = '  +0.79%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.42'
$ws.Range("E42").Value = '  +3.82%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.786'
$ws.Range("E43").Value = '  +1.04%  '

# Row 44
$ws.Range("D44").Value = '1.766.21'
$ws.Range("E44").Value = '  +1.67%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '92.58'
$ws.Range("E45").Value = '  +0.77%  '

# Row 46
$ws.Range("E46").Value = '  +2.56%  '

# Row 47
$ws.Range("E47").Value = '  +0.78%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0511'
$ws.Range("E48").Value = '  +1.04%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.410'
$ws.Range("E49").Value = '  +0.86%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.48'
$ws.Range("E50").Value = '  +3.74%  '

# Row 51
$ws.Range("E51").Value = '  +0.06%  '
